$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2375
$ws.Cells.Item(40, 10).Value = 2500
$ws.Cells.Item(40, 12).Value = 2500
$ws.Cells.Item(40, 14).Value = -2850
$ws.Cells.Item(62, 8).Value = 6130.615
$ws.Cells.Item(62, 9).Value = 4962.25
$ws.Cells.Item(62, 11).Value = 4962.25
$ws.Cells.Item(62, 13).Value = -4338.25
$ws.Cells.Item(64, 8).Value = 4400
$ws.Cells.Item(64, 9).Value = 5300
$ws.Cells.Item(64, 11).Value = 5300
$ws.Cells.Item(64, 13).Value = -5052
$ws.Cells.Item(65, 8).Value = 6130.615
$ws.Cells.Item(65, 9).Value = 4962.25
$ws.Cells.Item(65, 11).Value = 24811.25
$ws.Cells.Item(65, 13).Value = -21691.25
$ws.Cells.Item(67, 8).Value = 4400
$ws.Cells.Item(67, 9).Value = 5300
$ws.Cells.Item(67, 11).Value = 5300
$ws.Cells.Item(67, 13).Value = -4442
$ws.Cells.Item(76, 8).Value = 6292.1113
$ws.Cells.Item(76, 9).Value = 5176.2
$ws.Cells.Item(76, 11).Value = 5176.2
$ws.Cells.Item(76, 13).Value = -4861.2
$ws.Cells.Item(79, 8).Value = 6292.1113
$ws.Cells.Item(79, 9).Value = 5176.2
$ws.Cells.Item(79, 11).Value = 5176.2
$ws.Cells.Item(79, 13).Value = -4084.2
$ws.Cells.Item(116, 8).Value = 9299.444
$ws.Cells.Item(116, 9).Value = 9132.666999999999
$ws.Cells.Item(116, 11).Value = 9132.666999999999
$ws.Cells.Item(116, 13).Value = -5690.666999999999
$ws.Cells.Item(132, 8).Value = 1530.68
$ws.Cells.Item(132, 9).Value = 1563.65
$ws.Cells.Item(132, 11).Value = 4690.950000000001
$ws.Cells.Item(132, 13).Value = -2160.950000000001
$ws.Cells.Item(137, 8).Value = 3479.25
$ws.Cells.Item(137, 10).Value = 4656.7144
$ws.Cells.Item(137, 12).Value = 13970.1432
$ws.Cells.Item(137, 14).Value = -19070.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17424.441
$ws.Cells.Item(32, 9).Value = 8112.96
$ws.Cells.Item(32, 11).Value = 8112.96
$ws.Cells.Item(32, 13).Value = -7825.96
$ws.Cells.Item(74, 8).Value = 5299.6665
$ws.Cells.Item(74, 9).Value = 1415.2858
$ws.Cells.Item(74, 11).Value = 1415.2858
$ws.Cells.Item(74, 13).Value = -541.2858000000001
$ws.Cells.Item(77, 8).Value = 5299.6665
$ws.Cells.Item(77, 9).Value = 1415.2858
$ws.Cells.Item(77, 11).Value = 7076.429
$ws.Cells.Item(77, 13).Value = -2708.429
$ws.Cells.Item(119, 8).Value = 72000
$ws.Cells.Item(119, 10).Value = 72000
$ws.Cells.Item(119, 12).Value = 72000
$ws.Cells.Item(119, 14).Value = -81676
$ws.Cells.Item(132, 8).Value = 931.34485
$ws.Cells.Item(132, 9).Value = 947.7778
$ws.Cells.Item(132, 11).Value = 2843.3334
$ws.Cells.Item(132, 13).Value = -313.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1676.25
$ws.Cells.Item(20, 9).Value = 1568.6666
$ws.Cells.Item(20, 10).Value = 1999
$ws.Cells.Item(20, 11).Value = 1568.6666
$ws.Cells.Item(20, 12).Value = 1999
$ws.Cells.Item(20, 13).Value = -1321.6666
$ws.Cells.Item(20, 14).Value = -2493
$ws.Cells.Item(22, 8).Value = 536.82355
$ws.Cells.Item(22, 9).Value = 542.25
$ws.Cells.Item(22, 11).Value = 542.25
$ws.Cells.Item(22, 13).Value = -369.25
$ws.Cells.Item(105, 8).Value = 3678.9062
$ws.Cells.Item(105, 9).Value = 2991.35
$ws.Cells.Item(105, 11).Value = 2991.35
$ws.Cells.Item(105, 13).Value = -1244.35
$ws.Cells.Item(134, 8).Value = 2760.3
$ws.Cells.Item(134, 10).Value = 2615.75
$ws.Cells.Item(134, 12).Value = 7847.25
$ws.Cells.Item(134, 14).Value = -12917.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3644.125
$ws.Cells.Item(31, 9).Value = 2069.625
$ws.Cells.Item(31, 11).Value = 2069.625
$ws.Cells.Item(31, 13).Value = -1774.625
$ws.Cells.Item(34, 8).Value = 3644.125
$ws.Cells.Item(34, 9).Value = 2069.625
$ws.Cells.Item(34, 11).Value = 2069.625
$ws.Cells.Item(34, 13).Value = -1867.625
$ws.Cells.Item(122, 8).Value = 7635.1904
$ws.Cells.Item(122, 9).Value = 7356.643
$ws.Cells.Item(122, 11).Value = 22069.929
$ws.Cells.Item(122, 13).Value = -19619.929

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 77040070
$ws.Cells.Item(4, 10).Value = 7499.5
$ws.Cells.Item(4, 12).Value = 22498.5
$ws.Cells.Item(4, 14).Value = -22722.5
$ws.Cells.Item(5, 8).Value = 646.25
$ws.Cells.Item(5, 9).Value = 345.42856
$ws.Cells.Item(5, 10).Value = 1067.4
$ws.Cells.Item(5, 11).Value = 1036.28568
$ws.Cells.Item(5, 12).Value = 3202.2
$ws.Cells.Item(5, 13).Value = -924.28568
$ws.Cells.Item(5, 14).Value = -3426.2
$ws.Cells.Item(113, 8).Value = 1159.5714
$ws.Cells.Item(113, 9).Value = 1889.4
$ws.Cells.Item(113, 10).Value = 931.5
$ws.Cells.Item(113, 11).Value = 5668.200000000001
$ws.Cells.Item(113, 12).Value = 2794.5
$ws.Cells.Item(113, 13).Value = -3498.200000000001
$ws.Cells.Item(113, 14).Value = -7134.5
$ws.Cells.Item(132, 8).Value = 9849.75
$ws.Cells.Item(132, 10).Value = 9849.75
$ws.Cells.Item(132, 12).Value = 88647.75
$ws.Cells.Item(132, 14).Value = -93707.75
$ws.Cells.Item(135, 8).Value = 646.25
$ws.Cells.Item(135, 9).Value = 345.42856
$ws.Cells.Item(135, 10).Value = 1067.4
$ws.Cells.Item(135, 11).Value = 3108.85704
$ws.Cells.Item(135, 12).Value = 9606.6
$ws.Cells.Item(135, 13).Value = -573.8570399999999
$ws.Cells.Item(135, 14).Value = -14676.6
$ws.Cells.Item(140, 8).Value = 2578.3333
$ws.Cells.Item(140, 9).Value = 2578.3333
$ws.Cells.Item(140, 11).Value = 7734.999899999999
$ws.Cells.Item(140, 13).Value = -2554.999899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5978.8335
$ws.Cells.Item(70, 9).Value = 4976.5
$ws.Cells.Item(70, 10).Value = 6780.7
$ws.Cells.Item(70, 11).Value = 4976.5
$ws.Cells.Item(70, 12).Value = 6780.7
$ws.Cells.Item(70, 13).Value = -4706.5
$ws.Cells.Item(70, 14).Value = -7320.7
$ws.Cells.Item(73, 8).Value = 5978.8335
$ws.Cells.Item(73, 9).Value = 4976.5
$ws.Cells.Item(73, 10).Value = 6780.7
$ws.Cells.Item(73, 11).Value = 4976.5
$ws.Cells.Item(73, 12).Value = 6780.7
$ws.Cells.Item(73, 13).Value = -4040.5
$ws.Cells.Item(73, 14).Value = -8652.700000000001
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 922.6087
$ws.Cells.Item(107, 9).Value = 782.5
$ws.Cells.Item(107, 10).Value = 952.1053000000001
$ws.Cells.Item(107, 11).Value = 782.5
$ws.Cells.Item(107, 12).Value = 952.1053000000001
$ws.Cells.Item(107, 13).Value = 1137.5
$ws.Cells.Item(107, 14).Value = -4792.1053
$ws.Cells.Item(122, 8).Value = 1005068.4
$ws.Cells.Item(122, 9).Value = 171961.33
$ws.Cells.Item(122, 11).Value = 515883.99
$ws.Cells.Item(122, 13).Value = -513433.99

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2149.6667
$ws.Cells.Item(40, 9).Value = 2149.6667
$ws.Cells.Item(40, 11).Value = 2149.6667
$ws.Cells.Item(40, 13).Value = -2013.6667
$ws.Cells.Item(82, 8).Value = 1800.6666
$ws.Cells.Item(82, 9).Value = 2169
$ws.Cells.Item(82, 10).Value = 1465.8182
$ws.Cells.Item(82, 11).Value = 2169
$ws.Cells.Item(82, 12).Value = 1465.8182
$ws.Cells.Item(82, 13).Value = -1808
$ws.Cells.Item(82, 14).Value = -2187.8182
$ws.Cells.Item(85, 8).Value = 1800.6666
$ws.Cells.Item(85, 9).Value = 2169
$ws.Cells.Item(85, 10).Value = 1465.8182
$ws.Cells.Item(85, 11).Value = 2169
$ws.Cells.Item(85, 12).Value = 1465.8182
$ws.Cells.Item(85, 13).Value = -921
$ws.Cells.Item(85, 14).Value = -3961.8182
$ws.Cells.Item(136, 8).Value = 3076.7083
$ws.Cells.Item(136, 9).Value = 2944.8572
$ws.Cells.Item(136, 11).Value = 8834.571599999999
$ws.Cells.Item(136, 13).Value = -6284.571599999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1212.8
$ws.Cells.Item(107, 9).Value = 532.5714
$ws.Cells.Item(107, 11).Value = 1597.7142
$ws.Cells.Item(107, 13).Value = 322.2857999999999
$ws.Cells.Item(112, 8).Value = 45693.5
$ws.Cells.Item(112, 10).Value = 45693.5
$ws.Cells.Item(112, 12).Value = 45693.5
$ws.Cells.Item(112, 14).Value = -48647.5
$ws.Cells.Item(113, 8).Value = 1330.826
$ws.Cells.Item(113, 9).Value = 1336.1818
$ws.Cells.Item(113, 10).Value = 1325.9166
$ws.Cells.Item(113, 11).Value = 4008.5454
$ws.Cells.Item(113, 12).Value = 3977.7498
$ws.Cells.Item(113, 13).Value = -1838.5454
$ws.Cells.Item(113, 14).Value = -8317.7498
$ws.Cells.Item(138, 8).Value = 199984
$ws.Cells.Item(138, 10).Value = 199984
$ws.Cells.Item(138, 12).Value = 199984
$ws.Cells.Item(138, 14).Value = -210264
